$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "aggiunta click su ogni casella (MouseListener)"
# Fix typo in the backlog item text: "controllo errore" -> "controllo errori"
$ws.Range("C14").Value = "controllo errori"

# Mark the two remaining un-estimated backlog rows (controllo errori / caricamento
# vari Sudoku) with their estimated/true time, like the other rows already have
# in columns H (Tempo Stimato) and I (Tempo veritiero).
$ws.Range("H14").Value = "1gg"
$ws.Range("I14").Value = "1gg"
$ws.Range("H16").Value = "1gg"
$ws.Range("I16").Value = "1gg"

# Scroll the sheet so column B is at the left edge of the view.
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
